$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (F column) counts
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 99
$wsExpo.Range("F5").Value = 2652
$wsExpo.Range("F6").Value = 251
$wsExpo.Range("F7").Value = 383

# Sheet "全部类型" (all types) - same exhibitions repeated, update matching rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 99
$wsAll.Range("F5").Value = 2652
$wsAll.Range("F6").Value = 251
$wsAll.Range("F9").Value = 383
